$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "65.172.35"
$ws.Range("E2").Value = "  -0.12%  "

Set-TextValue "D3" "3.529.97"
$ws.Range("E3").Value = "  +3.62%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.10%  "

Set-TextValue "D5" "597.19"
$ws.Range("E5").Value = "  +2.57%  "

Set-TextValue "D6" "137.61"
$ws.Range("E6").Value = "  +1.04%  "

Set-TextValue "D7" "3.527.85"
$ws.Range("E7").Value = "  +3.56%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("E10").Value = "  +3.31%  "

Set-TextValue "D11" "6.89"
$ws.Range("E11").Value = "  -3.94%  "

$ws.Range("E12").Value = "  +3.35%  "

Set-TextValue "D13" "4.121.23"
$ws.Range("E13").Value = "  +3.59%  "

$ws.Range("E14").Value = "  +2.72%  "

$ws.Range("E15").Value = "  +5.11%  "

Set-TextValue "D16" "3.530.24"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("E17").Value = "  +1.39%  "

Set-TextValue "D18" "65.093.16"
$ws.Range("E18").Value = "  -0.11%  "

Set-TextValue "D19" "10.12"
$ws.Range("E19").Value = "  +4.24%  "

$ws.Range("E20").Value = "  +1.12%  "

Set-TextValue "D21" "14.21"

Set-TextValue "D22" "390.95"
$ws.Range("E22").Value = "  +2.31%  "

$ws.Range("E23").Value = "  +3.71%  "

Set-TextValue "D24" "3.668.44"
$ws.Range("E24").Value = "  +3.70%  "

Set-TextValue "D25" "73.44"
$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  +9.25%  "

Set-TextValue "D28" "7.78"
$ws.Range("E28").Value = "  +11.08%  "

$ws.Range("E29").Value = "  +0.00%  "

Set-TextValue "D30" "2.27"
$ws.Range("E30").Value = "  +3.50%  "

Set-TextValue "D31" "8.12"
$ws.Range("E31").Value = "  +1.62%  "

Set-TextValue "D32" "3.545.80"
$ws.Range("E32").Value = "  +4.02%  "

Set-TextValue "D34" "23.82"
$ws.Range("E34").Value = "  +5.20%  "

$ws.Range("E35").Value = "  +1.68%  "

$ws.Range("E36").Value = "  +16.67%  "

Set-TextValue "D37" "169.94"
$ws.Range("E37").Value = "  -0.32%  "

$ws.Range("E38").Value = "  +8.12%  "

Set-TextValue "D39" "6.90"
$ws.Range("E39").Value = "  +3.01%  "

Set-TextValue "D40" "4.97"
$ws.Range("E40").Value = "  +6.15%  "

Set-TextValue "D41" "0.0799"
$ws.Range("E41").Value = "  +6.03%  "

Set-TextValue "D42" "0.825"
$ws.Range("E42").Value = "  +1.52%  "

Set-TextValue "D43" "26.23"
$ws.Range("E43").Value = "  +18.62%  "

Set-TextValue "D44" "42.62"
$ws.Range("E44").Value = "  -1.69%  "

Set-TextValue "D45" "0.998"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D47" "1.19"
$ws.Range("E47").Value = "  +9.22%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D48" "1.67"
$ws.Range("E48").Value = "  +5.12%  "

$ws.Range("E49").Value = "  +4.82%  "

Set-TextValue "D50" "2.387.54"
$ws.Range("E50").Value = "  +10.15%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D51" "2.36"
$ws.Range("E51").Value = "  +20.18%  "
